# Shift every "Perioden" (year) label in column A forward by one year
# (2009 -> 2010, 2010 -> 2011, ... 2020 -> 2021). The "Perioden" column
# stores the year as TEXT (it round-trips through the shared-string
# table), so the new value has to be written back as text as well -
# otherwise Excel's normal type inference would store a numeric-looking
# string like "2010" as a number instead of text.
#
# Trick: prefixing the value with a leading single-quote forces Excel to
# treat it as text (this is the standard COM/VBA way of doing it), but it
# also stamps the cell with a "quote prefix" number format. Calling
# ClearFormats() right after removes that stray per-cell formatting again
# so the cell ends up with exactly the same (default) style it started
# with - only its text content changed, just like the source diff shows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

# Row 1 is the header ("Perioden" / "Postcode" / "SJV_GEMIDDELD") and is
# left untouched; the year values start on row 2.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $year = [int]$cell.Value2
    $newYear = $year + 1
    $cell.Value = "'" + [string]$newYear
    $cell.ClearFormats()
}
